# Apply the "create sensor_variable(s)" edit described by the commit:
#  - add two new worksheets: "Shelf Sensor" and "Sensor Data"
#  - add a new constant row (shelf_sensor_reg_size / INT) on "Constants"
#  - widen the "Constants" A column to fit the new, longer name
#  - change the selected range on "Shelf" (no content change)
#  - leave "Sensor Data" as the active sheet/tab when done

$wb = $excel.ActiveWorkbook

$constants = $wb.Worksheets.Item("Constants")
$shelf     = $wb.Worksheets.Item("Shelf")

# Grab a couple of already-fitted column widths from the existing sheets so
# the new sheets' bestfit-ish columns line up with the rest of the workbook
# instead of guessing arbitrary numbers.
$shelfColAWidth     = $shelf.Columns.Item(1).ColumnWidth      # narrow numeric/addr column
$constantsColAWidth = $constants.Columns.Item(1).ColumnWidth  # "variable_name"-ish column

# ---------------------------------------------------------------------
# 1. Constants: new row describing the sensor register block size
# ---------------------------------------------------------------------
$constants.Range("A4").Value = "shelf_sensor_reg_size"
$constants.Range("C4").Value = "INT"

# Column A needs to be a bit wider to fit "shelf_sensor_reg_size"
$constants.Columns.Item(1).ColumnWidth = 20

# ---------------------------------------------------------------------
# 2. Shelf: selection moved from a single cell to the header row
# ---------------------------------------------------------------------
$shelf.Activate()
[void]$shelf.Range("B1:E1").Select()

# ---------------------------------------------------------------------
# 3. New sheet "Shelf Sensor" (placed right after "Shelf")
# ---------------------------------------------------------------------
$shelfSensor = $wb.Worksheets.Add($null, $shelf)
$shelfSensor.Name = "Shelf Sensor"

$shelfSensor.Range("A1").Value = "base_addr"
$shelfSensor.Range("B1").Value = "variable_name"
$shelfSensor.Range("B1").HorizontalAlignment = -4108   # xlCenter

$shelfSensor.Range("A2").Value = 12000
$shelfSensor.Range("B2").Value  = "wPres0"
$shelfSensor.Range("B3").Value  = "wPres1"
$shelfSensor.Range("B4").Value  = "wPres2"
$shelfSensor.Range("B5").Value  = "wTemp0"
$shelfSensor.Range("B6").Value  = "wTemp1"
$shelfSensor.Range("B7").Value  = "wTemp2"
# filled bottom-up to match the original authoring order (aTemp, aRH, aCO2)
$shelfSensor.Range("B10").Value = "aTemp"
$shelfSensor.Range("B9").Value  = "aRH"
$shelfSensor.Range("B8").Value  = "aCO2"

$shelfSensor.Columns.Item(1).ColumnWidth = $shelfColAWidth
$shelfSensor.Columns.Item(2).ColumnWidth = $constantsColAWidth

$shelfSensor.Activate()
[void]$shelfSensor.Range("A3").Select()

# ---------------------------------------------------------------------
# 4. New sheet "Sensor Data" (placed right after "Shelf Sensor")
# ---------------------------------------------------------------------
$sensorData = $wb.Worksheets.Add($null, $shelfSensor)
$sensorData.Name = "Sensor Data"

$sensorData.Range("A1").Value = "variable_name"
$sensorData.Range("B1").Value = "type"
$sensorData.Range("C1").Value = "init_value"
$sensorData.Range("A1:C1").HorizontalAlignment = -4108  # xlCenter

$sensorRows = @("value", "state", "warn_u", "warn_l", "err_u", "err_l")
$r = 2
foreach ($name in $sensorRows) {
    $sensorData.Range("A$r").Value = $name
    $sensorData.Range("B$r").Value = "WORD"
    $sensorData.Range("C$r").Value = 0
    $r++
}

$sensorData.Columns.Item(1).ColumnWidth = $constantsColAWidth

# "Sensor Data" ends up being the active sheet/tab in the saved workbook
$sensorData.Activate()
[void]$sensorData.Range("A8").Select()
